$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '79.572.72'
$ws.Range('D3').Value = '3.204.32'
$ws.Range('E3').Value = '  +5.29%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '205.21'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.13%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '635.86'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.54%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +13.67%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.586'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +6.12%  '
$ws.Range('D10').Value = '3.202.26'
$ws.Range('E10').Value = '  +5.29%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.591'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +34.37%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.166'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +3.16%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.50'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +7.51%  '
$ws.Range('D14').Value = '3.790.05'
$ws.Range('E14').Value = '  +5.03%  '
$ws.Range('E15').Value = '  +18.70%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '31.85'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +7.74%  '
$ws.Range('D17').Value = '79.447.72'
$ws.Range('E17').Value = '  +4.04%  '
$ws.Range('D18').Value = '3.195.10'
$ws.Range('E18').Value = '  +5.09%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '14.52'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +7.54%  '
$ws.Range('E20').Value = '  +30.81%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '9.19'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.57%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '427.82'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +14.19%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.11'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +17.22%  '
$ws.Range('B24').Value = 'WrappedeETH'
$ws.Range('C24').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D24').Value = '3.368.23'
$ws.Range('E24').Value = '  +5.24%  '
$ws.Range('B25').Value = 'Aptos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '11.27'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +13.15%  '
$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '4.77'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +7.92%  '
$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '77.01'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.74%  '
$ws.Range('E28').Value = '  +0.17%  '
$ws.Range('E29').Value = '  +7.31%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '9.05'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +8.98%  '
$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.26%  '
$ws.Range('E32').Value = '  +5.04%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '524.08'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.30%  '
$ws.Range('E34').Value = '  +2.29%  '
$ws.Range('E35').Value = '  +27.88%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '22.95'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +9.67%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.120'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +12.93%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.999'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.408'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +5.78%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '164.86'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.06%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '192.23'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.28%  '
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.54'
$ws.Range('D44').ClearFormats()
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.824'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.13%  '
$ws.Range('E46').Value = '  +7.96%  '
$ws.Range('E47').Value = '  +2.87%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '43.12'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.38%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '25.84'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +14.97%  '
$ws.Range('E50').Value = '  +4.97%  '
$ws.Range('B51').Value = 'Filecoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '4.18'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +7.35%  '
